# Mixture_template.xlsx — "fixed data error for density"
#
# - B12/B13 on "Materials and Products" were mistakenly left at a
#   placeholder 100; correct them to 2 and 5 respectively.
# - Widen columns D/E/F on the same sheet so the longer header text fits.
# - Leave the cursor on B13 (the last cell touched) on "Materials and
#   Products", matching where the author's selection ended up.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Materials and Products")
$ws.Activate()

# --- fix the density data-entry error -------------------------------
$ws.Range("B12").Value = 2
$ws.Range("B13").Value = 5

# --- widen columns D, E and F ----------------------------------------
$ws.Columns.Item(4).ColumnWidth = 36.83333333333333
$ws.Columns.Item(5).ColumnWidth = 25.166666666666664
$ws.Columns.Item(6).ColumnWidth = 27.666666666666664

# --- leave the selection where the author left it ---------------------
$ws.Range("B13").Select()
